$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-14 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-15 Monday", 2) | Out-Null
$d.Content.Find.Execute("58×87=5046", $true, $false, $false, $false, $false, $true, 1, $false, "36×79=2844", 2) | Out-Null
$d.Content.Find.Execute("40×31=1240", $true, $false, $false, $false, $false, $true, 1, $false, "36×86=3096", 2) | Out-Null
$d.Content.Find.Execute("50×21=1050", $true, $false, $false, $false, $false, $true, 1, $false, "83×26=2158", 2) | Out-Null
$d.Content.Find.Execute("29×21=609", $true, $false, $false, $false, $false, $true, 1, $false, "30×71=2130", 2) | Out-Null
$d.Content.Find.Execute("93×33=3069", $true, $false, $false, $false, $false, $true, 1, $false, "78×77=6006", 2) | Out-Null
$d.Content.Find.Execute("31×60=1860", $true, $false, $false, $false, $false, $true, 1, $false, "46×80=3680", 2) | Out-Null
$d.Content.Find.Execute("91×96=8736", $true, $false, $false, $false, $false, $true, 1, $false, "72×62=4464", 2) | Out-Null
$d.Content.Find.Execute("40×58=2320", $true, $false, $false, $false, $false, $true, 1, $false, "79×19=1501", 2) | Out-Null
$d.Content.Find.Execute("82×24=1968", $true, $false, $false, $false, $false, $true, 1, $false, "78×43=3354", 2) | Out-Null
$d.Content.Find.Execute("41×86=3526", $true, $false, $false, $false, $false, $true, 1, $false, "41×36=1476", 2) | Out-Null
$d.Content.Find.Execute("86×23=1978", $true, $false, $false, $false, $false, $true, 1, $false, "57×19=1083", 2) | Out-Null
$d.Content.Find.Execute("89×40=3560", $true, $false, $false, $false, $false, $true, 1, $false, "68×74=5032", 2) | Out-Null
$d.Content.Find.Execute("90×76=6840", $true, $false, $false, $false, $false, $true, 1, $false, "84×70=5880", 2) | Out-Null
$d.Content.Find.Execute("91×47=4277", $true, $false, $false, $false, $false, $true, 1, $false, "77×11=847", 2) | Out-Null
$d.Content.Find.Execute("26×84=2184", $true, $false, $false, $false, $false, $true, 1, $false, "94×39=3666", 2) | Out-Null
$d.Content.Find.Execute("31×70=2170", $true, $false, $false, $false, $false, $true, 1, $false, "28×98=2744", 2) | Out-Null
$d.Content.Find.Execute("66×86=5676", $true, $false, $false, $false, $false, $true, 1, $false, "24×34=816", 2) | Out-Null
$d.Content.Find.Execute("97×20=1940", $true, $false, $false, $false, $false, $true, 1, $false, "99×28=2772", 2) | Out-Null
$d.Content.Find.Execute("75×67=5025", $true, $false, $false, $false, $false, $true, 1, $false, "85×73=6205", 2) | Out-Null
$d.Content.Find.Execute("41×51=2091", $true, $false, $false, $false, $false, $true, 1, $false, "60×90=5400", 2) | Out-Null
$d.Content.Find.Execute("82×33=2706", $true, $false, $false, $false, $false, $true, 1, $false, "34×48=1632", 2) | Out-Null
$d.Content.Find.Execute("37×68=2516", $true, $false, $false, $false, $false, $true, 1, $false, "66×80=5280", 2) | Out-Null
$d.Content.Find.Execute("32×66=2112", $true, $false, $false, $false, $false, $true, 1, $false, "29×48=1392", 2) | Out-Null
$d.Content.Find.Execute("86×99=8514", $true, $false, $false, $false, $false, $true, 1, $false, "42×41=1722", 2) | Out-Null
$d.Content.Find.Execute("48×24=1152", $true, $false, $false, $false, $false, $true, 1, $false, "13×66=858", 2) | Out-Null
